$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 264, shifting rows 264:339 down to 265:340.
$ws.Rows(264).Insert()

# Populate the newly inserted row 264 with the new record.
# Columns A,B,C,E,F,G,H,I,N,O,Q,R keep the same values the old row 264 had
# (Choclo / Dulce o Americano / Primera / $/malla 70 unidades / Región de
# Arica y Parinacota / 70 / Hortaliza); only D,J,K,L,M,P are new.
$ws.Range("A264").Value = 7
$ws.Range("B264").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C264").Value = "Ñuble"
$ws.Range("D264").Value = 45093
$ws.Range("E264").Value = 16
$ws.Range("F264").Value = 100112024
$ws.Range("G264").Value = "Choclo"
$ws.Range("H264").Value = "Dulce o Americano"
$ws.Range("I264").Value = "Primera"
$ws.Range("J264").Value = 30
$ws.Range("K264").Value = 16000
$ws.Range("L264").Value = 16000
$ws.Range("M264").Value = 16000
$ws.Range("N264").Value = "$/malla 70 unidades"
$ws.Range("O264").Value = "Región de Arica y Parinacota"
$ws.Range("P264").Value = 229
$ws.Range("Q264").Value = 70
$ws.Range("R264").Value = "Hortaliza"

# Keep the date column's display format consistent with the rest of column D.
$ws.Range("D264").NumberFormat = "YYYY-MM-DD HH:MM:SS"
